$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 50292.117
$ws.Range("I74").Value = 66376.234
$ws.Range("K74").Value = 66376.234
$ws.Range("M74").Value = -65440.234
$ws.Range("H77").Value = 50292.117
$ws.Range("I77").Value = 66376.234
$ws.Range("K77").Value = 331881.17
$ws.Range("M77").Value = -327201.17
$ws.Range("H137").Value = 2096.2
$ws.Range("I137").Value = 2555.75
$ws.Range("J137").Value = 1571
$ws.Range("K137").Value = 7667.25
$ws.Range("L137").Value = 4713
$ws.Range("M137").Value = -5117.25
$ws.Range("N137").Value = -9813

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 50000
$ws.Range("J23").Value = 50000
$ws.Range("L23").Value = 50000
$ws.Range("N23").Value = -50518
$ws.Range("H32").Value = 2972.611
$ws.Range("I32").Value = 2337.9375
$ws.Range("K32").Value = 2337.9375
$ws.Range("M32").Value = -2050.9375
$ws.Range("H122").Value = 168402.83
$ws.Range("I122").Value = 251298.75
$ws.Range("K122").Value = 753896.25
$ws.Range("M122").Value = -751446.25
$ws.Range("H132").Value = 1749.0682
$ws.Range("I132").Value = 1321.6571
$ws.Range("J132").Value = 3411.2222
$ws.Range("K132").Value = 3964.9713
$ws.Range("L132").Value = 10233.6666
$ws.Range("M132").Value = -1434.9713
$ws.Range("N132").Value = -15293.6666

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2115.5
$ws.Range("I107").Value = 2213.353
$ws.Range("J107").Value = 1699.625
$ws.Range("K107").Value = 2213.353
$ws.Range("L107").Value = 1699.625
$ws.Range("M107").Value = -293.3530000000001
$ws.Range("N107").Value = -5539.625

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 107969
$ws.Range("I31").Value = 168936.83
$ws.Range("J31").Value = 16517.25
$ws.Range("K31").Value = 168936.83
$ws.Range("L31").Value = 16517.25
$ws.Range("M31").Value = -168641.83
$ws.Range("N31").Value = -17107.25
$ws.Range("H34").Value = 107969
$ws.Range("I34").Value = 168936.83
$ws.Range("J34").Value = 16517.25
$ws.Range("K34").Value = 168936.83
$ws.Range("L34").Value = 16517.25
$ws.Range("M34").Value = -168734.83
$ws.Range("N34").Value = -16921.25
$ws.Range("H57").Value = 30237.5
$ws.Range("J57").Value = 30237.5
$ws.Range("L57").Value = 30237.5
$ws.Range("N57").Value = -31357.5
$ws.Range("H94").Value = 1425.6
$ws.Range("J94").Value = 1162.25
$ws.Range("L94").Value = 1162.25
$ws.Range("N94").Value = -2064.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 85
$ws.Range("I7").Value = 85
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 255
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -143
$ws.Range("N7").ClearContents()
$ws.Range("H8").Value = 1224
$ws.Range("I8").Value = 1224
$ws.Range("K8").Value = 3672
$ws.Range("M8").Value = -3533
$ws.Range("H76").Value = 1202.4
$ws.Range("I76").Value = 1202.4
$ws.Range("K76").Value = 3607.2
$ws.Range("M76").Value = -3224.2
$ws.Range("H79").Value = 1202.4
$ws.Range("I79").Value = 1202.4
$ws.Range("K79").Value = 3607.2
$ws.Range("M79").Value = -2281.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8915.294
$ws.Range("I70").Value = 7245.9
$ws.Range("J70").Value = 11300.143
$ws.Range("K70").Value = 7245.9
$ws.Range("L70").Value = 11300.143
$ws.Range("M70").Value = -6975.9
$ws.Range("N70").Value = -11840.143
$ws.Range("H73").Value = 8915.294
$ws.Range("I73").Value = 7245.9
$ws.Range("J73").Value = 11300.143
$ws.Range("K73").Value = 7245.9
$ws.Range("L73").Value = 11300.143
$ws.Range("M73").Value = -6309.9
$ws.Range("N73").Value = -13172.143
$ws.Range("H97").Value = 921.94116
$ws.Range("I97").Value = 788.2
$ws.Range("J97").Value = 1113
$ws.Range("K97").Value = 788.2
$ws.Range("L97").Value = 1113
$ws.Range("M97").Value = -292.2
$ws.Range("N97").Value = -2105
$ws.Range("H102").Value = 3719.9644
$ws.Range("I102").Value = 3772.182
$ws.Range("J102").Value = 3528.5
$ws.Range("K102").Value = 3772.182
$ws.Range("L102").Value = 3528.5
$ws.Range("M102").Value = -2150.182
$ws.Range("N102").Value = -6772.5
$ws.Range("H132").Value = 259089.31
$ws.Range("I132").Value = 315077.84
$ws.Range("K132").Value = 945233.52
$ws.Range("M132").Value = -942703.52

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 20836574
$ws.Range("I16").Value = 41668064
$ws.Range("J16").Value = 5083.25
$ws.Range("K16").Value = 41668064
$ws.Range("L16").Value = 5083.25
$ws.Range("M16").Value = -41667894
$ws.Range("N16").Value = -5423.25
$ws.Range("H22").Value = 710.375
$ws.Range("I22").Value = 414.66666
$ws.Range("J22").Value = 887.8
$ws.Range("K22").Value = 414.66666
$ws.Range("L22").Value = 887.8
$ws.Range("M22").Value = -119.66666
$ws.Range("N22").Value = -1477.8
$ws.Range("H27").Value = 710.375
$ws.Range("I27").Value = 414.66666
$ws.Range("J27").Value = 887.8
$ws.Range("K27").Value = 414.66666
$ws.Range("L27").Value = 887.8
$ws.Range("M27").Value = -307.66666
$ws.Range("N27").Value = -1101.8
$ws.Range("H43").Value = 21895
$ws.Range("J43").Value = 23158.334
$ws.Range("L43").Value = 23158.334
$ws.Range("N43").Value = -23544.334
$ws.Range("H68").Value = 10803.8125
$ws.Range("I68").Value = 3072.625
$ws.Range("J68").Value = 18535
$ws.Range("K68").Value = 3072.625
$ws.Range("L68").Value = 18535
$ws.Range("M68").Value = -2323.625
$ws.Range("N68").Value = -20033
$ws.Range("H71").Value = 10803.8125
$ws.Range("I71").Value = 3072.625
$ws.Range("J71").Value = 18535
$ws.Range("K71").Value = 15363.125
$ws.Range("L71").Value = 92675
$ws.Range("M71").Value = -11619.125
$ws.Range("N71").Value = -100163
$ws.Range("H100").Value = 1262119.2
$ws.Range("I100").Value = 1674659
$ws.Range("K100").Value = 1674659
$ws.Range("M100").Value = -1674118
$ws.Range("H101").Value = 13973
$ws.Range("J101").Value = 13973
$ws.Range("L101").Value = 13973
$ws.Range("N101").Value = -20463
$ws.Range("H122").Value = 4720.1113
$ws.Range("J122").Value = 5030
$ws.Range("L122").Value = 15090
$ws.Range("N122").Value = -19990
$ws.Range("H127").Value = 68607.5
$ws.Range("J127").Value = 68607.5
$ws.Range("L127").Value = 68607.5
$ws.Range("N127").Value = -78527.5
$ws.Range("H132").Value = 3905.3235
$ws.Range("I132").Value = 2876.0386
$ws.Range("J132").Value = 7250.5
$ws.Range("K132").Value = 8628.1158
$ws.Range("L132").Value = 21751.5
$ws.Range("M132").Value = -6098.1158
$ws.Range("N132").Value = -26811.5
$ws.Range("H140").Value = 45000
$ws.Range("J140").Value = 60000
$ws.Range("L140").Value = 60000
$ws.Range("N140").Value = -70360

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2384.4814
$ws.Range("I126").Value = 2409.8
$ws.Range("K126").Value = 7229.400000000001
$ws.Range("M126").Value = -4759.400000000001
$ws.Range("H132").Value = 3224.7097
$ws.Range("I132").Value = 3116.7144
$ws.Range("J132").Value = 4232.6665
$ws.Range("K132").Value = 9350.143199999999
$ws.Range("L132").Value = 12697.9995
$ws.Range("M132").Value = -6820.143199999999
$ws.Range("N132").Value = -17757.9995
